$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the error message text (shared string used by K5)
$ws.Range("K5").Value = "ERROR: La cantidad de semestres debe ser un numero entero"

# Update input values in row 3 (K3, L3)
$ws.Range("K3").Value = 11.4
$ws.Range("L3").Value = 12

# Update formulas in F5 / F6 (interest rate corrected, and F6 multiplier changed)
$ws.Range("F5").Formula = "= (0.3 * F2 * F3 * 0.0115) / (1 - (1 + 0.0115)^(-1 * (F3/2) * 12))"
$ws.Range("F6").Formula = "= (0.7 * F2 * F3 * 0.0115) / (1 - (1 + 0.0115)^(-1.5 * (F3/2) * 12))"

# Add new underlined cell M7 (new font/style), expanding used range to column M
$ws.Range("M7").Font.Underline = $true

# Update the active selection shown when the workbook is opened
[void]$ws.Range("E18").Select()
